$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "last updated" timestamp in A1 (21:16 -> 21:46)
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 21:46"

# Swap Huelva/Huesca rows (row 53 was Huesca/0, row 54 was Huelva/72;
# after the edit row 53 is Huelva/72 and row 54 is Huesca/0)
$ws.Range("A53").Value = "Huelva"
$ws.Range("C53").Value = 72

$ws.Range("A54").Value = "Huesca"
$ws.Range("C54").Value = 0
